$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "60.705.49"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.617.11"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.33%  "
Set-TextValue "D5" "514.96"
$ws.Range("E5").Value = "  +1.55%  "
Set-TextValue "D6" "154.45"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.12%  "
Set-TextValue "D8" "0.589"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "2.628.87"
$ws.Range("E9").Value = "  -0.14%  "
Set-TextValue "D10" "6.72"
$ws.Range("E10").Value = "  +4.17%  "
$ws.Range("E11").Value = "  -0.10%  "
Set-TextValue "D12" "0.347"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").Value = "3.073.93"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "60.647.96"
$ws.Range("E15").Value = "  +0.23%  "
Set-TextValue "D16" "21.74"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "2.622.14"
$ws.Range("E18").Value = "  +0.01%  "
Set-TextValue "D19" "4.76"
$ws.Range("E19").Value = "  -0.07%  "
Set-TextValue "D20" "357.67"
$ws.Range("E20").Value = "  +4.01%  "
Set-TextValue "D21" "10.68"
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  -0.03%  "
Set-TextValue "D24" "60.84"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "2.731.83"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("E27").Value = "  +0.82%  "
Set-TextValue "D28" "0.997"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "0.0₃0848"
$ws.Range("E29").Value = "  -0.78%  "
Set-TextValue "D30" "7.38"
$ws.Range("E30").Value = "  -2.24%  "
Set-TextValue "D31" "1.00"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E32").Value = "  +0.24%  "
Set-TextValue "D33" "5.99"
$ws.Range("E33").Value = "  +4.52%  "
$ws.Range("E34").Value = "  +1.13%  "
Set-TextValue "D35" "151.85"
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("E37").Value = "  -0.60%  "
Set-TextValue "D38" "0.885"
$ws.Range("E38").Value = "  +6.84%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D39" "0.856"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D40" "1.49"
$ws.Range("E40").Value = "  +1.02%  "
Set-TextValue "D41" "36.42"
$ws.Range("E41").Value = "  +2.49%  "
Set-TextValue "D42" "3.76"
$ws.Range("E42").Value = "  -0.42%  "
Set-TextValue "D43" "293.67"
$ws.Range("E43").Value = "  -4.59%  "
$ws.Range("E44").Value = "  +0.70%  "
Set-TextValue "D45" "0.623"
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("E46").Value = "  -2.28%  "
Set-TextValue "D47" "0.996"
$ws.Range("E47").Value = "  -0.38%  "
Set-TextValue "D48" "19.80"
$ws.Range("E48").Value = "  -0.07%  "
Set-TextValue "D49" "4.96"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("E50").Value = "  +0.30%  "
Set-TextValue "D51" "10.29"
$ws.Range("E51").Value = "  +0.26%  "
